$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL Results")

# --- Resize the saved window (best-effort; cosmetic window geometry) ---
$wb.Windows.Item(1).Width = 19440

# --- Change D8 to use the "quote prefix" text style while keeping its numeric value ---
# Build the quote-prefix style first via a scratch cell, then paste only the
# formatting onto D8 so the underlying value (2) is preserved as a number.
$ws.Range("A9:G9").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)
$ws.Range("A12").Value = "'Temp"
$ws.Range("A12").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("A12:G12").Delete()

# --- Add new row 10 (POGGIO A CAIANO) ---
$ws.Range("A9:G9").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)
$ws.Range("C10").Value = "POGGIO A CAIANO"
$ws.Range("F10").Value = "PO"
$ws.Range("A10").Value = "'100004"
$ws.Range("B10").Value = "G754"
$ws.Range("E10").Value = "'09"
$ws.Range("D10").Value = 100
$ws.Range("A10").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("G10").Value = 42919.777777777781

# --- Add new row 11 (NARDO') ---
$ws.Range("A9:G9").Copy()
$ws.Range("A11:G11").PasteSpecial(-4122)
$ws.Range("A11").Value = "'075052"
$ws.Range("C11").Value = "NARDO'"
$ws.Range("F11").Value = "LE"
$ws.Range("B11").Value = "F842"
$ws.Range("E11").Value = "'16"
$ws.Range("D11").Value = 75
$ws.Range("A11").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("G11").Value = 42926.736805555556

# --- Update the active selection to reflect where the author left off editing ---
$ws.Range("E13").Select()
